# Auto-generated: update cryptocurrency Price (D) and Volume(1h) (E) columns
# to reflect the refreshed data pull (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''29.304.47'
$ws.Range("E2").Value = '''  +0.18%  '
$ws.Range("D3").Value = '''1.859.63'
$ws.Range("E3").Value = '''  +0.05%  '
$ws.Range("E4").Value = '''  -0.02%  '
$ws.Range("D5").Value = '''0.7018'
$ws.Range("E5").Value = '''  +0.08%  '
$ws.Range("E6").Value = '''  -0.07%  '
$ws.Range("D7").Value = '''1.000'
$ws.Range("E7").Value = '''  +0.01%  '
$ws.Range("D8").Value = '''0.07793'
$ws.Range("E8").Value = '''  -3.63%  '
$ws.Range("D9").Value = '''0.3046'
$ws.Range("E9").Value = '''  +0.15%  '
$ws.Range("E10").Value = '''  +6.44%  '
$ws.Range("D11").Value = '''0.08144'
$ws.Range("E11").Value = '''  -0.45%  '
$ws.Range("D12").Value = '''1.871.45'
$ws.Range("E12").Value = '''  +0.53%  '
$ws.Range("D13").Value = '''5.211'
$ws.Range("E13").Value = '''  +0.57%  '
$ws.Range("D14").Value = '''0.7125'
$ws.Range("E14").Value = '''  -0.38%  '
$ws.Range("D15").Value = '''89.18'
$ws.Range("E15").Value = '''  -0.01%  '
$ws.Range("D16").Value = '''29.302.30'
$ws.Range("E16").Value = '''  +0.11%  '
$ws.Range("D17").Value = '''243.81'
$ws.Range("E17").Value = '''  +3.31%  '
$ws.Range("D18").Value = '''5.778'
$ws.Range("E18").Value = '''  +0.08%  '
$ws.Range("D19").Value = '''0.000007770'
$ws.Range("E19").Value = '''  -1.03%  '
$ws.Range("D20").Value = '''13.18'
$ws.Range("E20").Value = '''  -1.40%  '
$ws.Range("D21").Value = '''2.117.36'
$ws.Range("E21").Value = '''  +0.67%  '
$ws.Range("D22").Value = '''0.9999'
$ws.Range("E22").Value = '''  +0.09%  '
$ws.Range("E23").Value = '''  -0.03%  '
$ws.Range("D24").Value = '''7.525'
$ws.Range("E24").Value = '''  +0.97%  '
$ws.Range("D25").Value = '''162.16'
$ws.Range("E25").Value = '''  +0.20%  '
$ws.Range("D26").Value = '''8.860'
$ws.Range("E26").Value = '''  -1.30%  '
$ws.Range("D27").Value = '''0.1434'
$ws.Range("E27").Value = '''  -1.25%  '
$ws.Range("D28").Value = '''18.06'
$ws.Range("E28").Value = '''  +0.02%  '
$ws.Range("D29").Value = '''1.906'
$ws.Range("E29").Value = '''  -4.34%  '
$ws.Range("E30").Value = '''  -4.77%  '
$ws.Range("D31").Value = '''1.471'
$ws.Range("E31").Value = '''  -0.89%  '
$ws.Range("D32").Value = '''4.296'
$ws.Range("E32").Value = '''  -2.45%  '
$ws.Range("D33").Value = '''4.025'
$ws.Range("E33").Value = '''  -0.67%  '
$ws.Range("D34").Value = '''0.05154'
$ws.Range("E34").Value = '''  -0.98%  '
$ws.Range("D35").Value = '''1.180'
$ws.Range("E35").Value = '''  +0.84%  '
$ws.Range("D36").Value = '''0.7058'
$ws.Range("E36").Value = '''  -0.18%  '
$ws.Range("D37").Value = '''0.9971'
$ws.Range("E37").Value = '''  -0.41%  '
$ws.Range("E38").Value = '''  +0.50%  '
$ws.Range("D39").Value = '''0.01844'
$ws.Range("E39").Value = '''  -0.22%  '
$ws.Range("D40").Value = '''2.689'
$ws.Range("E40").Value = '''  -1.14%  '
$ws.Range("D41").Value = '''1.171.89'
$ws.Range("E41").Value = '''  +2.44%  '
$ws.Range("D42").Value = '''0.9124'
$ws.Range("E42").Value = '''  -1.23%  '
$ws.Range("D43").Value = '''5.998'
$ws.Range("E43").Value = '''  +0.50%  '
$ws.Range("D44").Value = '''71.24'
$ws.Range("E44").Value = '''  +0.54%  '
$ws.Range("D45").Value = '''0.4241'
$ws.Range("E45").Value = '''  -0.73%  '
$ws.Range("D46").Value = '''0.9999'
$ws.Range("E46").Value = '''  +0.05%  '
$ws.Range("D47").Value = '''101.23'
$ws.Range("E47").Value = '''  -1.99%  '
$ws.Range("D48").Value = '''0.5350'
$ws.Range("E48").Value = '''  -1.12%  '
$ws.Range("E49").Value = '''  -1.90%  '
$ws.Range("D50").Value = '''9.134'
$ws.Range("E50").Value = '''  -0.65%  '
$ws.Range("D51").Value = '''6.941'
$ws.Range("E51").Value = '''  +0.08%  '
